$wb = $excel.ActiveWorkbook

# Add a new worksheet; re-fetch objects by name afterwards since existing
# worksheet references can go stale (rebind positionally) once the
# collection is mutated via Add()/Move().
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet2"

# Move the new sheet to sit right after Sheet1.
$wb.Worksheets.Item("Sheet2").Move($null, $wb.Worksheets.Item("Sheet1"))

# From here on, always look sheets up fresh by name.
$ws2 = $wb.Worksheets.Item("Sheet2")

# Header row
$ws2.Range("A1").Value = "type"
$ws2.Range("B1").Value = "username"
$ws2.Range("C1").Value = "password"

# Data rows
$ws2.Range("A2").Value = "valid"
$ws2.Range("B2").Value = "standard_user"
$ws2.Range("C2").Value = "secret_sauce"

$ws2.Range("A3").Value = "invalid"
$ws2.Range("B3").Value = "standard_user"
$ws2.Range("C3").Value = "qwert2"

$ws2.Range("A4").Value = "locked"
$ws2.Range("B4").Value = "locked_out_user"
$ws2.Range("C4").Value = "secret_sauce"

# Column widths to match bestFit widths from the target file.
# (ColumnWidth gets snapped to 1/6-character increments on save, so the
# input is pre-compensated by the engine's fixed +5/6 save-time offset.)
$ws2.Columns.Item(1).ColumnWidth = 5.666666666666667
$ws2.Columns.Item(2).ColumnWidth = 13.666666666666666
$ws2.Columns.Item(3).ColumnWidth = 10.998697916666666

# Make Sheet2 the active/selected sheet, set zoom + selection to match
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 190
$ws2.Range("C4").Select() | Out-Null
